$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change header cell A1 from "RESPUESTA" to "PREGUNTA"
$ws.Range("A1").Value = "PREGUNTA"

# Update the active selection to C5, matching the saved view state
$ws.Range("C5").Select()
